$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AX9").Value = "Anders Forsberg, Alexander Hoffmann, David Isaksson"
$ws.Range("A13").Value = 131106312
$ws.Range("AB13").Value = "13:29"
$ws.Range("B13").Value = 92106
$ws.Range("E13").Value = 658
$ws.Range("F13").Value = "Rosenticka"
$ws.Range("G13").Value = "Fomitopsis rosea"
$ws.Range("H13").Value = "(Alb. & Schwein.:Fr.) P.Karst."
$ws.Range("Q13").Value = 601540
$ws.Range("R13").Value = 6992576
$ws.Range("X13").Value = "2025_0872"
$ws.Range("Z13").Value = "13:29"
$ws.Range("A14").Value = 131106314
$ws.Range("AB14").Value = "13:21"
$ws.Range("AX14").Value = "David Isaksson"
$ws.Range("B14").Value = 91808
$ws.Range("E14").Value = 1202
$ws.Range("F14").Value = "Ullticka"
$ws.Range("G14").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H14").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("I14").Value = "1"
$ws.Range("Q14").Value = 601556
$ws.Range("R14").Value = 6992605
$ws.Range("X14").Value = "2025_0870"
$ws.Range("Z14").Value = "13:21"
$ws.Range("A15").Value = 131106325
$ws.Range("AB15").Value = "12:21"
$ws.Range("AX15").Value = "Alexander Hoffmann"
$ws.Range("I15").Value = ""
$ws.Range("Q15").Value = 601615
$ws.Range("R15").Value = 6992785
$ws.Range("X15").Value = "2025_0858"
$ws.Range("Z15").Value = "12:21"
$ws.Range("AX18").Value = "Jennifer Lehikoinen, Linnéa Kjellberg"
$ws.Range("A20").Value = 131106310
$ws.Range("AB20").Value = "13:35"
$ws.Range("AX20").Value = "Alexander Hoffmann"
$ws.Range("I20").Value = ""
$ws.Range("J20").Value = ""
$ws.Range("Q20").Value = 601470
$ws.Range("R20").Value = 6992568
$ws.Range("X20").Value = "2025_0874"
$ws.Range("Z20").Value = "13:35"
$ws.Range("A21").Value = 131106313
$ws.Range("AB21").Value = "13:29"
$ws.Range("AX21").Value = "David Isaksson"
$ws.Range("I21").Value = "1"
$ws.Range("J21").Value = "mycel"
$ws.Range("Q21").Value = 601530
$ws.Range("R21").Value = 6992589
$ws.Range("X21").Value = "2025_0871"
$ws.Range("Z21").Value = "13:29"
$ws.Range("A24").Value = 131106311
$ws.Range("AB24").Value = "13:32"
$ws.Range("AX24").Value = "Alexander Hoffmann"
$ws.Range("J24").Value = ""
$ws.Range("Q24").Value = 601498
$ws.Range("R24").Value = 6992583
$ws.Range("X24").Value = "2025_0873"
$ws.Range("Z24").Value = "13:32"
$ws.Range("A25").Value = 131106327
$ws.Range("AB25").Value = "12:10"
$ws.Range("AX25").Value = "Alexander Hoffmann, David Isaksson"
$ws.Range("J25").Value = "mycel"
$ws.Range("Q25").Value = 601607
$ws.Range("R25").Value = 6992789
$ws.Range("X25").Value = "2025_0856"
$ws.Range("Z25").Value = "12:10"
$ws.Range("A27").Value = 131106330
$ws.Range("AB27").Value = "12:06"
$ws.Range("AX27").Value = "Alexander Hoffmann, David Isaksson"
$ws.Range("B27").Value = 92106
$ws.Range("D27").Value = "NT"
$ws.Range("E27").Value = 658
$ws.Range("F27").Value = "Rosenticka"
$ws.Range("G27").Value = "Fomitopsis rosea"
$ws.Range("H27").Value = "(Alb. & Schwein.:Fr.) P.Karst."
$ws.Range("J27").Value = ""
$ws.Range("Q27").Value = 601607
$ws.Range("R27").Value = 6992782
$ws.Range("X27").Value = "2025_0853"
$ws.Range("Z27").Value = "12:06"
$ws.Range("A28").Value = 131106329
$ws.Range("AB28").Value = "12:09"
$ws.Range("AX28").Value = "Alexander Hoffmann, David Isaksson"
$ws.Range("B28").Value = 92267
$ws.Range("D28").Value = "VU"
$ws.Range("E28").Value = 1209
$ws.Range("F28").Value = "Rynkskinn"
$ws.Range("G28").Value = "Hermanssonia centrifuga"
$ws.Range("H28").Value = "(P. Karst.) Zmitr."
$ws.Range("J28").Value = "mycel"
$ws.Range("Q28").Value = 601609
$ws.Range("R28").Value = 6992789
$ws.Range("X28").Value = "2025_0854"
$ws.Range("Z28").Value = "12:09"
